$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Date column (B2:B4): was stored as literal text strings like
# "2025-01-02"; replace with real Excel date serials formatted mm-dd-yy,
# left-aligned (matches style used elsewhere in the sheet).
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = 45659

$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = 45660

$ws.Range("B4").HorizontalAlignment = -4131
$ws.Range("B4").NumberFormat = "mm-dd-yy"
$ws.Range("B4").Value = 45661

# --- Row 4 previously only had A4/B4 populated (attendance record for
# 2025-01-04 was incomplete). Add an empty time-formatted In-Time cell
# (C4) and an actual Out-Time value (D4, 18:34 stored as a time serial).
$ws.Range("C4").NumberFormat = "h:mm"

$ws.Range("D4").HorizontalAlignment = -4131
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("D4").Value = 0.77361111111111114

# --- Selection moves to the newly-populated D4 cell.
$ws.Range("D4").Select()
